$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new rows before row 18 so the old rows 18-23 (summary row,
# average row, blank rows) shift down to rows 20-25, matching the target.
# ---------------------------------------------------------------------------
$ws.Rows("18:19").Insert()

# ---------------------------------------------------------------------------
# New row 18 -> Example 20
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = 20
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "V"
$ws.Range("D18").Value = "X"
$ws.Range("E18").Value = 4

# ---------------------------------------------------------------------------
# New row 19 -> Example 21
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = 21
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = "V"
$ws.Range("D19").Value = "X"
$ws.Range("E19").Value = 6
$ws.Range("F19").Formula = "=230/234"
$ws.Range("G19").Formula = "=486/476"
$ws.Range("H19").Formula = "=321/271"
$ws.Range("L19").Value = 0.76529999999999998

# ---------------------------------------------------------------------------
# New text entries, written in the same order the target sharedStrings
# table introduces them: Cons. Match Ratio, the Italian VTM-compliance
# note, "Wrong choice (T instead of C)", then the two MV-pair labels.
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = "Cons. Match Ratio"
$ws.Range("N18").Value = "In questo caso, il fatto che il costruttore di candidati non sia VTM-compliant fa sì che la scelta sia quella corretta. Per capire meglio, vedi il vettore dei SAD per questo esempio."
$ws.Range("M19").Value = "Wrong choice (T instead of C)"
$ws.Range("I19").Value = "(inf;180°)"
$ws.Range("J19").Value = "(inf;180°)"
$ws.Range("K19").Value = "(1;90°)"

# ---------------------------------------------------------------------------
# Row 21 (old row 19): AVERAGE grows to include the two new rows, and a new
# "Constr. Match" ratio (COUNTIF) appears in column D.
# ---------------------------------------------------------------------------
$ws.Range("B21").Formula = "=AVERAGE(B2:B19)"
$ws.Range("D21").Formula = '=COUNTIF(D2:D19,"V")/(COUNTIF(D2:D19,"V")+COUNTIF(D2:D19,"X"))'

# ---------------------------------------------------------------------------
# Style fix-ups: writing .Value / .Formula above reset number formats on
# the touched cells, so restore the correct look. Cloning an existing
# cell's format (copy / paste-special) preserves the exact style entry,
# including the quotePrefix flag, rather than approximating it.
# ---------------------------------------------------------------------------

# L18 / L19: quote-prefixed plain-number look, same as L15/L16/L17.
$ws.Range("L15").Copy() | Out-Null
$ws.Range("L18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("L15").Copy() | Out-Null
$ws.Range("L19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# D21: quote-prefixed 0.00 format, right aligned (new style - build before
# the L19 number-format tweak below so it is registered first, matching
# the target style order).
$ws.Range("L15").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D21").NumberFormat = "0.00"
$ws.Range("D21").HorizontalAlignment = -4152      # xlRight

# L19: quote-prefixed 0.000 format (new style, registered after D21's).
$ws.Range("L19").NumberFormat = "0.000"

Write-Output "done"
